$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3126.5312  # ALC!H76: 3088.0444 -> 3126.5312
$ws.Cells.Item(76, 9).Value = 3001.75  # ALC!I76: 2999.36 -> 3001.75
$ws.Cells.Item(76, 10).Value = 3201.4  # ALC!J76: 3198.9 -> 3201.4
$ws.Cells.Item(76, 11).Value = 3001.75  # ALC!K76: 2999.36 -> 3001.75
$ws.Cells.Item(76, 12).Value = 3201.4  # ALC!L76: 3198.9 -> 3201.4
$ws.Cells.Item(76, 13).Value = -2686.75  # ALC!M76: -2684.36 -> -2686.75
$ws.Cells.Item(76, 14).Value = -3831.4  # ALC!N76: -3828.9 -> -3831.4

$ws.Cells.Item(79, 8).Value = 3126.5312  # ALC!H79: 3088.0444 -> 3126.5312
$ws.Cells.Item(79, 9).Value = 3001.75  # ALC!I79: 2999.36 -> 3001.75
$ws.Cells.Item(79, 10).Value = 3201.4  # ALC!J79: 3198.9 -> 3201.4
$ws.Cells.Item(79, 11).Value = 3001.75  # ALC!K79: 2999.36 -> 3001.75
$ws.Cells.Item(79, 12).Value = 3201.4  # ALC!L79: 3198.9 -> 3201.4
$ws.Cells.Item(79, 13).Value = -1909.75  # ALC!M79: -1907.36 -> -1909.75
$ws.Cells.Item(79, 14).Value = -5385.4  # ALC!N79: -5382.9 -> -5385.4

$ws.Cells.Item(92, 8).Value = 792.1429000000001  # ALC!H92: 376.47058 -> 792.1429000000001
$ws.Cells.Item(92, 9).Value = 590.8333  # ALC!I92: 292.85715 -> 590.8333
$ws.Cells.Item(92, 10).Value = 2000  # ALC!J92: 766.6667 -> 2000
$ws.Cells.Item(92, 11).Value = 590.8333  # ALC!K92: 292.85715 -> 590.8333
$ws.Cells.Item(92, 12).Value = 2000  # ALC!L92: 766.6667 -> 2000
$ws.Cells.Item(92, 13).Value = 657.1667  # ALC!M92: 955.14285 -> 657.1667
$ws.Cells.Item(92, 14).Value = -4496  # ALC!N92: -3262.6667 -> -4496

$ws.Cells.Item(101, 8).Value = 1147.1666  # ALC!H101: 1059.8334 -> 1147.1666
$ws.Cells.Item(101, 9).Value = 725.75  # ALC!I101: 594.75 -> 725.75
$ws.Cells.Item(101, 11).Value = 2177.25  # ALC!K101: 1784.25 -> 2177.25
$ws.Cells.Item(101, 13).Value = -555.25  # ALC!M101: -162.25 -> -555.25

$ws.Cells.Item(106, 8).Value = 40845.1  # ALC!H106: 33073.6 -> 40845.1
$ws.Cells.Item(106, 9).Value = 1569.4445  # ALC!I106: 1697.1538 -> 1569.4445
$ws.Cells.Item(106, 10).Value = 72979.73  # ALC!J106: 67064.75 -> 72979.73
$ws.Cells.Item(106, 11).Value = 1569.4445  # ALC!K106: 1697.1538 -> 1569.4445
$ws.Cells.Item(106, 12).Value = 72979.73  # ALC!L106: 67064.75 -> 72979.73
$ws.Cells.Item(106, 13).Value = -938.4445000000001  # ALC!M106: -1066.1538 -> -938.4445000000001
$ws.Cells.Item(106, 14).Value = -74241.73  # ALC!N106: -68326.75 -> -74241.73

$ws.Cells.Item(129, 8).Value = 1458.9524  # ALC!H129: 1501.9487 -> 1458.9524
$ws.Cells.Item(129, 10).Value = 1448.4242  # ALC!J129: 1503.2667 -> 1448.4242
$ws.Cells.Item(129, 12).Value = 4345.2726  # ALC!L129: 4509.800099999999 -> 4345.2726
$ws.Cells.Item(129, 14).Value = -14345.2726  # ALC!N129: -14509.8001 -> -14345.2726

$ws.Cells.Item(134, 8).Value = 48277.777  # ALC!H134: 46900 -> 48277.777
$ws.Cells.Item(134, 10).Value = 48277.777  # ALC!J134: 46900 -> 48277.777
$ws.Cells.Item(134, 12).Value = 48277.777  # ALC!L134: 46900 -> 48277.777
$ws.Cells.Item(134, 14).Value = -58417.777  # ALC!N134: -57040 -> -58417.777

$ws.Cells.Item(137, 8).Value = 3525.698  # ALC!H137: 3185.2834 -> 3525.698
$ws.Cells.Item(137, 9).Value = 1068.3143  # ALC!I137: 1001.1539 -> 1068.3143
$ws.Cells.Item(137, 10).Value = 8303.944  # ALC!J137: 7241.524 -> 8303.944
$ws.Cells.Item(137, 11).Value = 3204.9429  # ALC!K137: 3003.4617 -> 3204.9429
$ws.Cells.Item(137, 12).Value = 24911.832  # ALC!L137: 21724.572 -> 24911.832
$ws.Cells.Item(137, 13).Value = -654.9429  # ALC!M137: -453.4616999999998 -> -654.9429
$ws.Cells.Item(137, 14).Value = -30011.832  # ALC!N137: -26824.572 -> -30011.832

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1448.6774  # ARM!H61: 1602 -> 1448.6774
$ws.Cells.Item(61, 9).Value = 1086.3636  # ARM!I61: 1252.8125 -> 1086.3636
$ws.Cells.Item(61, 10).Value = 2334.3333  # ARM!J61: 2109.9092 -> 2334.3333
$ws.Cells.Item(61, 11).Value = 1086.3636  # ARM!K61: 1252.8125 -> 1086.3636
$ws.Cells.Item(61, 12).Value = 2334.3333  # ARM!L61: 2109.9092 -> 2334.3333
$ws.Cells.Item(61, 13).Value = -874.3635999999999  # ARM!M61: -1040.8125 -> -874.3635999999999
$ws.Cells.Item(61, 14).Value = -2758.3333  # ARM!N61: -2533.9092 -> -2758.3333

$ws.Cells.Item(74, 8).Value = 1390.9807  # ARM!H74: 1650.625 -> 1390.9807
$ws.Cells.Item(74, 9).Value = 1378.921  # ARM!I74: 1583.5161 -> 1378.921
$ws.Cells.Item(74, 10).Value = 1423.7142  # ARM!J74: 1881.7778 -> 1423.7142
$ws.Cells.Item(74, 11).Value = 1378.921  # ARM!K74: 1583.5161 -> 1378.921
$ws.Cells.Item(74, 12).Value = 1423.7142  # ARM!L74: 1881.7778 -> 1423.7142
$ws.Cells.Item(74, 13).Value = -504.921  # ARM!M74: -709.5161000000001 -> -504.921
$ws.Cells.Item(74, 14).Value = -3171.7142  # ARM!N74: -3629.7778 -> -3171.7142

$ws.Cells.Item(77, 8).Value = 1390.9807  # ARM!H77: 1650.625 -> 1390.9807
$ws.Cells.Item(77, 9).Value = 1378.921  # ARM!I77: 1583.5161 -> 1378.921
$ws.Cells.Item(77, 10).Value = 1423.7142  # ARM!J77: 1881.7778 -> 1423.7142
$ws.Cells.Item(77, 11).Value = 6894.605  # ARM!K77: 7917.5805 -> 6894.605
$ws.Cells.Item(77, 12).Value = 7118.571  # ARM!L77: 9408.889000000001 -> 7118.571
$ws.Cells.Item(77, 13).Value = -2526.605  # ARM!M77: -3549.5805 -> -2526.605
$ws.Cells.Item(77, 14).Value = -15854.571  # ARM!N77: -18144.889 -> -15854.571

$ws.Cells.Item(86, 10).Value = 0  # ARM!J86: 28000 -> 0
$ws.Cells.Item(86, 12).Value = 0  # ARM!L86: 28000 -> 0
$ws.Cells.Item(86, 14).ClearContents()  # ARM!N86: -30372 -> (removed)

$ws.Cells.Item(89, 10).Value = 0  # ARM!J89: 28000 -> 0
$ws.Cells.Item(89, 12).Value = 0  # ARM!L89: 84000 -> 0
$ws.Cells.Item(89, 14).ClearContents()  # ARM!N89: -95856 -> (removed)

$ws.Cells.Item(102, 8).Value = 22407  # ARM!H102: 35261.668 -> 22407
$ws.Cells.Item(102, 10).Value = 36061.668  # ARM!J102: 101935 -> 36061.668
$ws.Cells.Item(102, 12).Value = 36061.668  # ARM!L102: 101935 -> 36061.668
$ws.Cells.Item(102, 14).Value = -39305.668  # ARM!N102: -105179 -> -39305.668

$ws.Cells.Item(136, 8).Value = 1448.6774  # ARM!H136: 1602 -> 1448.6774
$ws.Cells.Item(136, 9).Value = 1086.3636  # ARM!I136: 1252.8125 -> 1086.3636
$ws.Cells.Item(136, 10).Value = 2334.3333  # ARM!J136: 2109.9092 -> 2334.3333
$ws.Cells.Item(136, 11).Value = 3259.0908  # ARM!K136: 3758.4375 -> 3259.0908
$ws.Cells.Item(136, 12).Value = 7002.999899999999  # ARM!L136: 6329.7276 -> 7002.999899999999
$ws.Cells.Item(136, 13).Value = -709.0907999999999  # ARM!M136: -1208.4375 -> -709.0907999999999
$ws.Cells.Item(136, 14).Value = -12102.9999  # ARM!N136: -11429.7276 -> -12102.9999

$ws.Cells.Item(139, 8).Value = 49924.75  # ARM!H139: 48539.8 -> 49924.75
$ws.Cells.Item(139, 10).Value = 49924.75  # ARM!J139: 48539.8 -> 49924.75
$ws.Cells.Item(139, 12).Value = 49924.75  # ARM!L139: 48539.8 -> 49924.75
$ws.Cells.Item(139, 14).Value = -60204.75  # ARM!N139: -58819.8 -> -60204.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1817.2667  # BSM!H94: 1430 -> 1817.2667
$ws.Cells.Item(94, 9).Value = 1811.3572  # BSM!I94: 1430 -> 1811.3572
$ws.Cells.Item(94, 10).Value = 1900  # BSM!J94: 0 -> 1900
$ws.Cells.Item(94, 11).Value = 1811.3572  # BSM!K94: 1430 -> 1811.3572
$ws.Cells.Item(94, 12).Value = 1900  # BSM!L94: 0 -> 1900
$ws.Cells.Item(94, 13).Value = -1360.3572  # BSM!M94: -979 -> -1360.3572
$ws.Cells.Item(94, 14).Value = -2802  # BSM!N94: None -> -2802

$ws.Cells.Item(99, 8).Value = 1906.9445  # BSM!H99: 1813.125 -> 1906.9445
$ws.Cells.Item(99, 9).Value = 1845.3572  # BSM!I99: 1666.9231 -> 1845.3572
$ws.Cells.Item(99, 10).Value = 2122.5  # BSM!J99: 2446.6667 -> 2122.5
$ws.Cells.Item(99, 11).Value = 1845.3572  # BSM!K99: 1666.9231 -> 1845.3572
$ws.Cells.Item(99, 12).Value = 2122.5  # BSM!L99: 2446.6667 -> 2122.5
$ws.Cells.Item(99, 13).Value = -347.3571999999999  # BSM!M99: -168.9231 -> -347.3571999999999
$ws.Cells.Item(99, 14).Value = -5118.5  # BSM!N99: -5442.6667 -> -5118.5

$ws.Cells.Item(105, 8).Value = 3015.0334  # BSM!H105: 3342.2173 -> 3015.0334
$ws.Cells.Item(105, 9).Value = 2358.2354  # BSM!I105: 2552.9167 -> 2358.2354
$ws.Cells.Item(105, 10).Value = 3873.923  # BSM!J105: 4203.273 -> 3873.923
$ws.Cells.Item(105, 11).Value = 2358.2354  # BSM!K105: 2552.9167 -> 2358.2354
$ws.Cells.Item(105, 12).Value = 3873.923  # BSM!L105: 4203.273 -> 3873.923
$ws.Cells.Item(105, 13).Value = -611.2354  # BSM!M105: -805.9167000000002 -> -611.2354
$ws.Cells.Item(105, 14).Value = -7367.923  # BSM!N105: -7697.273 -> -7367.923

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3053.33  # CRP!H31: 3164.84 -> 3053.33
$ws.Cells.Item(31, 9).Value = 896.7143  # CRP!I31: 929.55554 -> 896.7143
$ws.Cells.Item(31, 10).Value = 3626.6077  # CRP!J31: 3655.5122 -> 3626.6077
$ws.Cells.Item(31, 11).Value = 896.7143  # CRP!K31: 929.55554 -> 896.7143
$ws.Cells.Item(31, 12).Value = 3626.6077  # CRP!L31: 3655.5122 -> 3626.6077
$ws.Cells.Item(31, 13).Value = -601.7143  # CRP!M31: -634.55554 -> -601.7143
$ws.Cells.Item(31, 14).Value = -4216.6077  # CRP!N31: -4245.5122 -> -4216.6077

$ws.Cells.Item(34, 8).Value = 3053.33  # CRP!H34: 3164.84 -> 3053.33
$ws.Cells.Item(34, 9).Value = 896.7143  # CRP!I34: 929.55554 -> 896.7143
$ws.Cells.Item(34, 10).Value = 3626.6077  # CRP!J34: 3655.5122 -> 3626.6077
$ws.Cells.Item(34, 11).Value = 896.7143  # CRP!K34: 929.55554 -> 896.7143
$ws.Cells.Item(34, 12).Value = 3626.6077  # CRP!L34: 3655.5122 -> 3626.6077
$ws.Cells.Item(34, 13).Value = -694.7143  # CRP!M34: -727.55554 -> -694.7143
$ws.Cells.Item(34, 14).Value = -4030.6077  # CRP!N34: -4059.5122 -> -4030.6077

$ws.Cells.Item(105, 8).Value = 5111  # CRP!H105: 2548.524 -> 5111
$ws.Cells.Item(105, 9).Value = 4749.875  # CRP!I105: 2514.0588 -> 4749.875
$ws.Cells.Item(105, 10).Value = 8000  # CRP!J105: 2695 -> 8000
$ws.Cells.Item(105, 11).Value = 4749.875  # CRP!K105: 2514.0588 -> 4749.875
$ws.Cells.Item(105, 12).Value = 8000  # CRP!L105: 2695 -> 8000
$ws.Cells.Item(105, 13).Value = -3002.875  # CRP!M105: -767.0587999999998 -> -3002.875
$ws.Cells.Item(105, 14).Value = -11494  # CRP!N105: -6189 -> -11494

$ws.Cells.Item(132, 8).Value = 43099.97  # CRP!H132: 56153.42 -> 43099.97
$ws.Cells.Item(132, 9).Value = 1517.5652  # CRP!I132: 1821.8889 -> 1517.5652
$ws.Cells.Item(132, 10).Value = 130045  # CRP!J132: 178399.38 -> 130045
$ws.Cells.Item(132, 11).Value = 4552.6956  # CRP!K132: 5465.6667 -> 4552.6956
$ws.Cells.Item(132, 12).Value = 390135  # CRP!L132: 535198.14 -> 390135
$ws.Cells.Item(132, 13).Value = -2022.6956  # CRP!M132: -2935.6667 -> -2022.6956
$ws.Cells.Item(132, 14).Value = -395195  # CRP!N132: -540258.14 -> -395195

$ws.Cells.Item(135, 8).Value = 57890  # CRP!H135: 63096 -> 57890
$ws.Cells.Item(135, 10).Value = 57890  # CRP!J135: 63096 -> 57890
$ws.Cells.Item(135, 12).Value = 57890  # CRP!L135: 63096 -> 57890
$ws.Cells.Item(135, 14).Value = -68030  # CRP!N135: -73236 -> -68030

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 356.82144  # CUL!H12: 369.18518 -> 356.82144
$ws.Cells.Item(12, 9).Value = 222.5  # CUL!I12: 183.9 -> 222.5
$ws.Cells.Item(12, 10).Value = 410.55  # CUL!J12: 478.17648 -> 410.55
$ws.Cells.Item(12, 11).Value = 667.5  # CUL!K12: 551.7 -> 667.5
$ws.Cells.Item(12, 12).Value = 1231.65  # CUL!L12: 1434.52944 -> 1231.65
$ws.Cells.Item(12, 13).Value = -494.5  # CUL!M12: -378.7 -> -494.5
$ws.Cells.Item(12, 14).Value = -1577.65  # CUL!N12: -1780.52944 -> -1577.65

$ws.Cells.Item(109, 8).Value = 2170.16  # CUL!H109: 2196.8333 -> 2170.16
$ws.Cells.Item(109, 9).Value = 1255.6923  # CUL!I109: 1286.4615 -> 1255.6923
$ws.Cells.Item(109, 10).Value = 3160.8333  # CUL!J109: 3272.7273 -> 3160.8333
$ws.Cells.Item(109, 11).Value = 3767.0769  # CUL!K109: 3859.3845 -> 3767.0769
$ws.Cells.Item(109, 12).Value = 9482.499899999999  # CUL!L109: 9818.1819 -> 9482.499899999999
$ws.Cells.Item(109, 13).Value = -2727.0769  # CUL!M109: -2819.3845 -> -2727.0769
$ws.Cells.Item(109, 14).Value = -11562.4999  # CUL!N109: -11898.1819 -> -11562.4999

$ws.Cells.Item(117, 8).Value = 966.6667  # CUL!H117: 1212.5 -> 966.6667
$ws.Cells.Item(117, 9).Value = 633.3333  # CUL!I117: 1500 -> 633.3333
$ws.Cells.Item(117, 10).Value = 1300  # CUL!J117: 925 -> 1300
$ws.Cells.Item(117, 11).Value = 1899.9999  # CUL!K117: 4500 -> 1899.9999
$ws.Cells.Item(117, 12).Value = 3900  # CUL!L117: 2775 -> 3900
$ws.Cells.Item(117, 13).Value = 1542.0001  # CUL!M117: -1058 -> 1542.0001
$ws.Cells.Item(117, 14).Value = -10784  # CUL!N117: -9659 -> -10784

$ws.Cells.Item(122, 8).Value = 4527.4136  # CUL!H122: 4140.4688 -> 4527.4136
$ws.Cells.Item(122, 9).Value = 788.3333  # CUL!I122: 749.5 -> 788.3333
$ws.Cells.Item(122, 11).Value = 7094.9997  # CUL!K122: 6745.5 -> 7094.9997
$ws.Cells.Item(122, 13).Value = -4644.9997  # CUL!M122: -4295.5 -> -4644.9997

$ws.Cells.Item(131, 8).Value = 3004.2075  # CUL!H131: 3452.25 -> 3004.2075
$ws.Cells.Item(131, 9).Value = 7140.4  # CUL!I131: 8173.5386 -> 7140.4
$ws.Cells.Item(131, 10).Value = 1371.5  # CUL!J131: 1472.3549 -> 1371.5
$ws.Cells.Item(131, 11).Value = 21421.2  # CUL!K131: 24520.6158 -> 21421.2
$ws.Cells.Item(131, 12).Value = 4114.5  # CUL!L131: 4417.0647 -> 4114.5
$ws.Cells.Item(131, 13).Value = -16381.2  # CUL!M131: -19480.6158 -> -16381.2
$ws.Cells.Item(131, 14).Value = -14194.5  # CUL!N131: -14497.0647 -> -14194.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 42154  # GSM!H104: 41210.168 -> 42154
$ws.Cells.Item(104, 10).Value = 42154  # GSM!J104: 41210.168 -> 42154
$ws.Cells.Item(104, 12).Value = 42154  # GSM!L104: 41210.168 -> 42154
$ws.Cells.Item(104, 14).Value = -49142  # GSM!N104: -48198.168 -> -49142

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 899.0323  # LTW!H93: 1131.7812 -> 899.0323
$ws.Cells.Item(93, 9).Value = 298.69565  # LTW!I93: 634.0417 -> 298.69565
$ws.Cells.Item(93, 11).Value = 298.69565  # LTW!K93: 634.0417 -> 298.69565
$ws.Cells.Item(93, 13).Value = 949.30435  # LTW!M93: 613.9583 -> 949.30435

$ws.Cells.Item(104, 8).Value = 0  # LTW!H104: 32546.666 -> 0
$ws.Cells.Item(104, 10).Value = 0  # LTW!J104: 32546.666 -> 0
$ws.Cells.Item(104, 12).Value = 0  # LTW!L104: 32546.666 -> 0
$ws.Cells.Item(104, 14).ClearContents()  # LTW!N104: -39534.666 -> (removed)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 539.6875  # WVR!H100: 502.05554 -> 539.6875
$ws.Cells.Item(100, 9).Value = 336.66666  # WVR!I100: 311.33334 -> 336.66666
$ws.Cells.Item(100, 10).Value = 661.5  # WVR!J100: 692.7778 -> 661.5
$ws.Cells.Item(100, 11).Value = 673.33332  # WVR!K100: 622.66668 -> 673.33332
$ws.Cells.Item(100, 12).Value = 1323  # WVR!L100: 1385.5556 -> 1323
$ws.Cells.Item(100, 13).Value = -132.33332  # WVR!M100: -81.66668000000004 -> -132.33332
$ws.Cells.Item(100, 14).Value = -2405  # WVR!N100: -2467.5556 -> -2405

$ws.Cells.Item(103, 8).Value = 45720.5  # WVR!H103: 47062.668 -> 45720.5
$ws.Cells.Item(103, 10).Value = 45720.5  # WVR!J103: 47062.668 -> 45720.5
$ws.Cells.Item(103, 12).Value = 45720.5  # WVR!L103: 47062.668 -> 45720.5
$ws.Cells.Item(103, 14).Value = -48064.5  # WVR!N103: -49406.668 -> -48064.5

$ws.Cells.Item(104, 8).Value = 44980.57  # WVR!H104: 46145.332 -> 44980.57
$ws.Cells.Item(104, 10).Value = 44980.57  # WVR!J104: 46145.332 -> 44980.57
$ws.Cells.Item(104, 12).Value = 44980.57  # WVR!L104: 46145.332 -> 44980.57
$ws.Cells.Item(104, 14).Value = -51968.57  # WVR!N104: -53133.332 -> -51968.57

$ws.Cells.Item(135, 8).Value = 33505  # WVR!H135: 29378.75 -> 33505
$ws.Cells.Item(135, 10).Value = 33505  # WVR!J135: 29378.75 -> 33505
$ws.Cells.Item(135, 12).Value = 33505  # WVR!L135: 29378.75 -> 33505
$ws.Cells.Item(135, 14).Value = -43645  # WVR!N135: -39518.75 -> -43645
